$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("TC_ID")

# Update header row
$ws.Range("A1").Value = "TC_ID"
$ws.Range("B1").Value = "ExecutinFlag"

# Update data rows with new test case names
$ws.Range("A2").Value = "001_GoodSignin"
$ws.Range("B2").Value = "Y"

$ws.Range("A3").Value = "002_BadSignin"
$ws.Range("B3").Value = "Y"

# Remove old row 4 (Validate_login_complete_test_3 / Y)
$ws.Rows.Item(4).Delete()
